$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one more year column (R), mirroring the formatting
# of the existing last column (Q) for rows 3 through 34.
$ws.Range("Q3:Q34").Copy($ws.Range("R3:R34"))

# Row 4 holds the year headers.
$ws.Range("R4").Value = 2021

# Data rows: set the new 2021 figures (column R) for each indicator row.
$ws.Range("R5").Value = 109
$ws.Range("R6").Value = 74
$ws.Range("R7").Value = 35
$ws.Range("R8").Value = 36
$ws.Range("R9").Value = 35
$ws.Range("R10").Value = 1
$ws.Range("R11").Value = 15
$ws.Range("R12").Value = 8
$ws.Range("R13").Value = 7
$ws.Range("R14").Value = 12
$ws.Range("R15").Value = 7
$ws.Range("R16").Value = 5
$ws.Range("R17").Value = "-"
$ws.Range("R18").Value = "-"
$ws.Range("R19").Value = "-"
$ws.Range("R20").Value = 17
$ws.Range("R21").Value = 8
$ws.Range("R22").Value = 9
$ws.Range("R23").Value = 9
$ws.Range("R24").Value = 7
$ws.Range("R25").Value = 2
$ws.Range("R26").Value = 20
$ws.Range("R27").Value = 9
$ws.Range("R28").Value = 11
$ws.Range("R29").Value = "-"
$ws.Range("R30").Value = "-"
$ws.Range("R31").Value = "-"
$ws.Range("R32").Value = "-"
$ws.Range("R33").Value = "-"
$ws.Range("R34").Value = "-"

# Match the final cursor/selection position recorded in the workbook.
$ws.Range("R35").Select()
